$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 10
$ws.Range("A10").Value = "Lavet vejledning og procedure til rename og move"
$ws.Range("B10").Value = "Tool Specialist"
$ws.Range("C10").Value = 43887
$ws.Range("D10").Value = 0.39583333333333331
$ws.Range("E10").Value = 0.4375
$ws.Range("F10").Value = 0.0625
$ws.Range("F10").NumberFormat = "h:mm"

# Row 11
$ws.Range("A11").Value = "Færdiggør vejledning og procedure til rename og move"
$ws.Range("B11").Value = "Tool Specialist"
$ws.Range("C11").Value = 43887
$ws.Range("D11").Value = 0.55208333333333337
$ws.Range("E11").Value = 0.5625
$ws.Range("F11").Value = 0.0069444444444444441
$ws.Range("F11").NumberFormat = "h:mm"

# Row 12
$ws.Range("A12").Value = "Lav AD04"
$ws.Range("B12").Value = "System Analyst "
$ws.Range("C12").Value = 43887
$ws.Range("D12").Value = 0.5625
$ws.Range("E12").Value = 0.58333333333333337
$ws.Range("F12").Value = 0.0069444444444444441
$ws.Range("F12").NumberFormat = "h:mm"

# Row 13
$ws.Range("A13").Value = "Review DOM05"
$ws.Range("B13").Value = "Reviewer"
$ws.Range("C13").Value = 43887
$ws.Range("D13").Value = 0.58333333333333337
$ws.Range("E13").Value = 0.58680555555555558
$ws.Range("F13").Value = 0.0069444444444444441
$ws.Range("F13").NumberFormat = "h:mm"

# Row 14
$ws.Range("A14").Value = "Rettet og færdiggjort ENV08 efter review"
$ws.Range("B14").Value = "Tool Specialist"
$ws.Range("C14").Value = 43887
$ws.Range("D14").Value = 0.61111111111111105
$ws.Range("E14").Value = 0.625
$ws.Range("F14").Value = 0.0069444444444444441
$ws.Range("F14").NumberFormat = "h:mm"

# Update selection to match diff
$ws.Range("D15").Select()
